$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 1878696.1
$ws.Range("I64").Value = 3207424.2
$ws.Range("J64").Value = 2844.7058
$ws.Range("K64").Value = 3207424.2
$ws.Range("L64").Value = 2844.7058
$ws.Range("M64").Value = -3207176.2
$ws.Range("N64").Value = -3340.7058

$ws.Range("H67").Value = 1878696.1
$ws.Range("I67").Value = 3207424.2
$ws.Range("J67").Value = 2844.7058
$ws.Range("K67").Value = 3207424.2
$ws.Range("L67").Value = 2844.7058
$ws.Range("M67").Value = -3206566.2
$ws.Range("N67").Value = -4560.7058

$ws.Range("H129").Value = 1220.2903
$ws.Range("I129").Value = 379.15384
$ws.Range("J129").Value = 1827.7778
$ws.Range("K129").Value = 1137.46152
$ws.Range("L129").Value = 5483.3334
$ws.Range("M129").Value = 3862.53848
$ws.Range("N129").Value = -15483.3334

$ws.Range("H138").Value = 9261606
$ws.Range("J138").Value = 2527.8572
$ws.Range("L138").Value = 7583.571599999999
$ws.Range("N138").Value = -17863.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9394
$ws.Range("I32").Value = 9831.508
$ws.Range("J32").Value = 7814.1113
$ws.Range("K32").Value = 9831.508
$ws.Range("L32").Value = 7814.1113
$ws.Range("M32").Value = -9544.508
$ws.Range("N32").Value = -8388.1113

$ws.Range("H122").Value = 6141.3335
$ws.Range("I122").Value = 7283
$ws.Range("J122").Value = 2716.3333
$ws.Range("K122").Value = 21849
$ws.Range("L122").Value = 8148.999899999999
$ws.Range("M122").Value = -19399
$ws.Range("N122").Value = -13048.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 47391.43
$ws.Range("J50").Value = 47391.43
$ws.Range("L50").Value = 47391.43
$ws.Range("N50").Value = -48539.43

$ws.Range("H94").Value = 534.1316
$ws.Range("I94").Value = 437.65518
$ws.Range("J94").Value = 845
$ws.Range("K94").Value = 437.65518
$ws.Range("L94").Value = 845
$ws.Range("M94").Value = 13.34482000000003
$ws.Range("N94").Value = -1747

$ws.Range("H99").Value = 1055
$ws.Range("I99").Value = 1060
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 1060
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = 438
$ws.Range("N99").Value = -3996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 3984.261
$ws.Range("I94").Value = 3013.75
$ws.Range("J94").Value = 4501.8667
$ws.Range("K94").Value = 3013.75
$ws.Range("L94").Value = 4501.8667
$ws.Range("M94").Value = -2562.75
$ws.Range("N94").Value = -5403.8667

$ws.Range("H122").Value = 4307.3335
$ws.Range("I122").Value = 4307.3335
$ws.Range("K122").Value = 12922.0005
$ws.Range("M122").Value = -10472.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6057.756
$ws.Range("I3").Value = 3927.7273
$ws.Range("J3").Value = 6838.7666
$ws.Range("K3").Value = 11783.1819
$ws.Range("L3").Value = 20516.2998
$ws.Range("M3").Value = -11671.1819
$ws.Range("N3").Value = -20740.2998

$ws.Range("H5").Value = 653.16
$ws.Range("I5").Value = 274.94736
$ws.Range("J5").Value = 1850.8334
$ws.Range("K5").Value = 824.84208
$ws.Range("L5").Value = 5552.5002
$ws.Range("M5").Value = -712.84208
$ws.Range("N5").Value = -5776.5002

$ws.Range("H122").Value = 1245.3077
$ws.Range("I122").Value = 1702
$ws.Range("J122").Value = 712.5
$ws.Range("K122").Value = 15318
$ws.Range("L122").Value = 6412.5
$ws.Range("M122").Value = -12868
$ws.Range("N122").Value = -11312.5

$ws.Range("H126").Value = 5500
$ws.Range("J126").Value = 6750
$ws.Range("L126").Value = 20250
$ws.Range("N126").Value = -30130

$ws.Range("H135").Value = 653.16
$ws.Range("I135").Value = 274.94736
$ws.Range("J135").Value = 1850.8334
$ws.Range("K135").Value = 2474.52624
$ws.Range("L135").Value = 16657.5006
$ws.Range("M135").Value = 60.47375999999986
$ws.Range("N135").Value = -21727.5006

$ws.Range("H136").Value = 2792.9412
$ws.Range("I136").Value = 1497.1428
$ws.Range("J136").Value = 3700
$ws.Range("K136").Value = 4491.428400000001
$ws.Range("L136").Value = 11100
$ws.Range("M136").Value = 608.5715999999993
$ws.Range("N136").Value = -21300

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8705.387000000001
$ws.Range("I70").Value = 12193.25
$ws.Range("J70").Value = 4519.95
$ws.Range("K70").Value = 12193.25
$ws.Range("L70").Value = 4519.95
$ws.Range("M70").Value = -11923.25
$ws.Range("N70").Value = -5059.95

$ws.Range("H73").Value = 8705.387000000001
$ws.Range("I73").Value = 12193.25
$ws.Range("J73").Value = 4519.95
$ws.Range("K73").Value = 12193.25
$ws.Range("L73").Value = 4519.95
$ws.Range("M73").Value = -11257.25
$ws.Range("N73").Value = -6391.95

$ws.Range("H122").Value = 2900665.2
$ws.Range("I122").Value = 4446093.5
$ws.Range("J122").Value = 2987.25
$ws.Range("K122").Value = 13338280.5
$ws.Range("L122").Value = 8961.75
$ws.Range("M122").Value = -13335830.5
$ws.Range("N122").Value = -13861.75

$ws.Range("H126").Value = 3548.1724
$ws.Range("I126").Value = 1932.4286
$ws.Range("J126").Value = 5056.2
$ws.Range("K126").Value = 5797.2858
$ws.Range("L126").Value = 15168.6
$ws.Range("M126").Value = -3327.2858
$ws.Range("N126").Value = -20108.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4765.3726
$ws.Range("I7").Value = 4735.8276
$ws.Range("J7").Value = 4804.3184
$ws.Range("K7").Value = 4735.8276
$ws.Range("L7").Value = 4804.3184
$ws.Range("M7").Value = -4623.8276
$ws.Range("N7").Value = -5028.3184

$ws.Range("H126").Value = 4765.3726
$ws.Range("I126").Value = 4735.8276
$ws.Range("J126").Value = 4804.3184
$ws.Range("K126").Value = 14207.4828
$ws.Range("L126").Value = 14412.9552
$ws.Range("M126").Value = -11737.4828
$ws.Range("N126").Value = -19352.9552

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 35000
$ws.Range("J63").Value = 35000
$ws.Range("L63").Value = 35000
$ws.Range("N63").Value = -36248

$ws.Range("H66").Value = 35000
$ws.Range("J66").Value = 35000
$ws.Range("L66").Value = 105000
$ws.Range("N66").Value = -111240

$ws.Range("H81").Value = 970.1111
$ws.Range("I81").Value = 632.2222
$ws.Range("J81").Value = 1308
$ws.Range("K81").Value = 1264.4444
$ws.Range("L81").Value = 2616
$ws.Range("M81").Value = -203.4444000000001
$ws.Range("N81").Value = -4738

$ws.Range("H84").Value = 970.1111
$ws.Range("I84").Value = 632.2222
$ws.Range("J84").Value = 1308
$ws.Range("K84").Value = 6322.222000000001
$ws.Range("L84").Value = 13080
$ws.Range("M84").Value = -1018.222000000001
$ws.Range("N84").Value = -23688

$ws.Range("H132").Value = 1177.3768
$ws.Range("I132").Value = 1094.305
$ws.Range("J132").Value = 1667.5
$ws.Range("K132").Value = 3282.915
$ws.Range("L132").Value = 5002.5
$ws.Range("M132").Value = -752.915
$ws.Range("N132").Value = -10062.5

$ws.Range("H136").Value = 1130.4688
$ws.Range("I136").Value = 833.10345
$ws.Range("K136").Value = 2499.31035
$ws.Range("M136").Value = 50.68965000000026
